$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the sample Pokemon rows (bulbasaur/charmander/ponyta) with the
# new set (Charizard/Eevee/Pikachu) and their ids/abilities.

# Row 2: Charizard
$ws.Range("B2").Value = "Charizard"
# Row 3: Eevee - new ability text
$ws.Range("C3").Value = "run-away, adaptability, anticipation"
# Row 3: Eevee name
$ws.Range("B3").Value = "Eevee"
# Row 4: Pikachu
$ws.Range("B4").Value = "Pikachu"
$ws.Range("C4").Value = "static, lightning-rod"
# Row 2: abilities for Charizard
$ws.Range("C2").Value = "blaze, solar-power"

# Update the numeric ids
$ws.Range("A2").Value = 6
$ws.Range("A3").Value = 133
$ws.Range("A4").Value = 25

# Update the selection to C4 (as reflected in the sheetView selection diff)
$ws.Range("C4").Select()
